$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in A5:A12 while keeping their existing style/formatting
$ws.Range("A5:A12").ClearContents()

# Update the selection to match the new active range
$ws.Range("A5:A12").Select()
